# "Generate Report for Handback"
#
# The handback round-trip for both target locales (zh-cn, de-de) has
# completed: the status text flips from "Ready for handoff" to
# "Handed back: in sync with en-US" everywhere it appears (Overview +
# per-locale sheets), and each locale sheet's two data rows get their
# "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns populated (they were blank placeholders before).
# A couple of report columns also get widened so the new, longer
# content is readable.

$wb = $excel.ActiveWorkbook

$mdUrl50ba = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/865e966a73debd1ba21973ece1e28af3a0d96289/e2e/50ba612b-7526-4296-a382-71bb777d8ff3.md"
$mdUrlC5e9 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/865e966a73debd1ba21973ece1e28af3a0d96289/e2e/c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md"

# ---------------------------------------------------------------------
# 1. Global status text: every "Ready for handoff" cell (Overview!E2:F3
#    and the Status column on both locale sheets) now reads as handed
#    back / in sync.
# ---------------------------------------------------------------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("Ready for handoff", "Handed back: in sync with en-US") | Out-Null
}

# ---------------------------------------------------------------------
# 2. zh-cn sheet: fill in the handback columns for both rows and widen
#    the Status / Latest Target File / Latest Handback File columns.
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl50ba, [Type]::Missing, [Type]::Missing, "50ba612b-7526-4296-a382-71bb777d8ff3.md")
$wsZh.Range("I2").Value = "50ba612b-7526-4296-a382-71bb777d8ff3.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl50ba, [Type]::Missing, [Type]::Missing, "50ba612b-7526-4296-a382-71bb777d8ff3.md")
$wsZh.Range("J2").Value = "50ba612b-7526-4296-a382-71bb777d8ff3.1ac77371c04877361c8ba536dda1ca67f9f0ba0d.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-20 09:51:02"

$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrlC5e9, [Type]::Missing, [Type]::Missing, "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md")
$wsZh.Range("I3").Value = "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrlC5e9, [Type]::Missing, [Type]::Missing, "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md")
$wsZh.Range("J3").Value = "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.ad48ab32875541a5093fbaf5e6146cd0f835a526.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-10-20 09:51:02"

$wsZh.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZh.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsZh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape, its own filenames/timestamp.
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl50ba, [Type]::Missing, [Type]::Missing, "50ba612b-7526-4296-a382-71bb777d8ff3.md")
$wsDe.Range("I2").Value = "50ba612b-7526-4296-a382-71bb777d8ff3.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl50ba, [Type]::Missing, [Type]::Missing, "50ba612b-7526-4296-a382-71bb777d8ff3.md")
$wsDe.Range("J2").Value = "50ba612b-7526-4296-a382-71bb777d8ff3.1ac77371c04877361c8ba536dda1ca67f9f0ba0d.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-20 09:51:20"

$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrlC5e9, [Type]::Missing, [Type]::Missing, "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md")
$wsDe.Range("I3").Value = "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrlC5e9, [Type]::Missing, [Type]::Missing, "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.md")
$wsDe.Range("J3").Value = "c5e9cb3d-b39e-487f-a343-dc6aacf129d4.ad48ab32875541a5093fbaf5e6146cd0f835a526.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-20 09:51:20"

$wsDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDe.Columns.Item(9).ColumnWidth = 39.166666666666664
$wsDe.Columns.Item(10).ColumnWidth = 39.166666666666664

# ---------------------------------------------------------------------
# 4. Overview sheet: widen the two locale-status columns (E, F) to
#    match the longer "Handed back: in sync with en-US" text.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668
